# Update the date in the title paragraph
$d = $word.ActiveDocument
$d.Content.Find.Execute("2026-01-16 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-17 Saturday", 2) | Out-Null

# Update the division problems in the table, addressed by (row, column)
# to avoid ambiguity since some old/new values collide across cells.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "83÷6=13, 5"  # was "79÷5=15, 4"
$t.Cell(1, 2).Range.Text = "31÷3=10, 1"  # was "95÷3=31, 2"
$t.Cell(1, 3).Range.Text = "16÷7=2, 2"  # was "80÷8=10, 0"
$t.Cell(1, 4).Range.Text = "88÷4=22, 0"  # was "31÷5=6, 1"
$t.Cell(1, 5).Range.Text = "54÷4=13, 2"  # was "22÷4=5, 2"
$t.Cell(5, 1).Range.Text = "55÷4=13, 3"  # was "32÷3=10, 2"
$t.Cell(5, 2).Range.Text = "33÷4=8, 1"  # was "23÷2=11, 1"
$t.Cell(5, 3).Range.Text = "34÷6=5, 4"  # was "43÷3=14, 1"
$t.Cell(5, 4).Range.Text = "71÷5=14, 1"  # was "84÷4=21, 0"
$t.Cell(5, 5).Range.Text = "46÷9=5, 1"  # was "29÷4=7, 1"
$t.Cell(9, 1).Range.Text = "82÷8=10, 2"  # was "60÷4=15, 0"
$t.Cell(9, 2).Range.Text = "71÷4=17, 3"  # was "51÷8=6, 3"
$t.Cell(9, 3).Range.Text = "38÷5=7, 3"  # was "93÷7=13, 2"
$t.Cell(9, 4).Range.Text = "19÷8=2, 3"  # was "45÷7=6, 3"
$t.Cell(9, 5).Range.Text = "92÷2=46, 0"  # was "99÷8=12, 3"
$t.Cell(13, 1).Range.Text = "16÷6=2, 4"  # was "43÷5=8, 3"
$t.Cell(13, 2).Range.Text = "71÷7=10, 1"  # was "30÷7=4, 2"
$t.Cell(13, 3).Range.Text = "46÷5=9, 1"  # was "86÷2=43, 0"
$t.Cell(13, 4).Range.Text = "85÷5=17, 0"  # was "59÷9=6, 5"
$t.Cell(13, 5).Range.Text = "15÷3=5, 0"  # was "44÷7=6, 2"
$t.Cell(17, 1).Range.Text = "23÷7=3, 2"  # was "87÷6=14, 3"
$t.Cell(17, 2).Range.Text = "10÷2=5, 0"  # was "81÷3=27, 0"
$t.Cell(17, 3).Range.Text = "32÷3=10, 2"  # was "92÷2=46, 0"
$t.Cell(17, 4).Range.Text = "64÷9=7, 1"  # was "70÷2=35, 0"
$t.Cell(17, 5).Range.Text = "20÷2=10, 0"  # was "16÷5=3, 1"
